$wb = $excel.ActiveWorkbook

# ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5304.25
$ws.Range("J40").Value = 6331.6665
$ws.Range("L40").Value = 6331.6665
$ws.Range("N40").Value = -6681.6665
$ws.Range("H111").Value = 870.75
$ws.Range("I111").Value = 744.5
$ws.Range("J111").Value = 997
$ws.Range("K111").Value = 2233.5
$ws.Range("L111").Value = 2991
$ws.Range("M111").Value = 833.5
$ws.Range("N111").Value = -9125
$ws.Range("H125").Value = 257074
$ws.Range("I125").Value = 4500
$ws.Range("K125").Value = 40500
$ws.Range("M125").Value = -38040
$ws.Range("H137").Value = 2362.0833
$ws.Range("I137").Value = 1913.2
$ws.Range("K137").Value = 5739.6
$ws.Range("M137").Value = -3189.6

# ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 577.6875
$ws.Range("I2").Value = 511.26666
$ws.Range("K2").Value = 511.26666
$ws.Range("M2").Value = -398.26666
$ws.Range("H74").Value = 898
$ws.Range("I74").Value = 875
$ws.Range("J74").Value = 913.3333
$ws.Range("K74").Value = 875
$ws.Range("L74").Value = 913.3333
$ws.Range("M74").Value = -1
$ws.Range("N74").Value = -2661.3333
$ws.Range("H77").Value = 898
$ws.Range("I77").Value = 875
$ws.Range("J77").Value = 913.3333
$ws.Range("K77").Value = 4375
$ws.Range("L77").Value = 4566.6665
$ws.Range("M77").Value = -7
$ws.Range("N77").Value = -13302.6665
$ws.Range("H88").Value = 2317.4614
$ws.Range("J88").Value = 3225.2222
$ws.Range("L88").Value = 3225.2222
$ws.Range("N88").Value = -4037.2222
$ws.Range("H91").Value = 2317.4614
$ws.Range("J91").Value = 3225.2222
$ws.Range("L91").Value = 3225.2222
$ws.Range("N91").Value = -6033.2222
$ws.Range("H116").Value = 577.6875
$ws.Range("I116").Value = 511.26666
$ws.Range("K116").Value = 511.26666
$ws.Range("M116").Value = 1782.73334
$ws.Range("H132").Value = 2949.6667
$ws.Range("J132").Value = 3950
$ws.Range("L132").Value = 11850
$ws.Range("N132").Value = -16910

# BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 577.6875
$ws.Range("I3").Value = 511.26666
$ws.Range("K3").Value = 511.26666
$ws.Range("M3").Value = -397.26666
$ws.Range("H7").Value = 333681.66
$ws.Range("I7").Value = 500022.5
$ws.Range("J7").Value = 1000
$ws.Range("K7").Value = 500022.5
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = -499909.5
$ws.Range("N7").Value = -1226

# CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1998.2354
$ws.Range("I31").Value = 1620.4166
$ws.Range("J31").Value = 2905
$ws.Range("K31").Value = 1620.4166
$ws.Range("L31").Value = 2905
$ws.Range("M31").Value = -1325.4166
$ws.Range("N31").Value = -3495
$ws.Range("H34").Value = 1998.2354
$ws.Range("I34").Value = 1620.4166
$ws.Range("J34").Value = 2905
$ws.Range("K34").Value = 1620.4166
$ws.Range("L34").Value = 2905
$ws.Range("M34").Value = -1418.4166
$ws.Range("N34").Value = -3309
$ws.Range("H51").Value = 20000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H61").Value = 20000
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H132").Value = 3097.077
$ws.Range("I132").Value = 1920.25
$ws.Range("J132").Value = 4980
$ws.Range("K132").Value = 5760.75
$ws.Range("L132").Value = 14940
$ws.Range("M132").Value = -3230.75
$ws.Range("N132").Value = -20000

# CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 121.85714
$ws.Range("I12").Value = 50.5
$ws.Range("J12").Value = 150.4
$ws.Range("K12").Value = 151.5
$ws.Range("L12").Value = 451.2
$ws.Range("M12").Value = 21.5
$ws.Range("N12").Value = -797.2
$ws.Range("H34").Value = 1491.6666
$ws.Range("I34").Value = 450
$ws.Range("J34").Value = 1586.3636
$ws.Range("K34").Value = 1350
$ws.Range("L34").Value = 4759.0908
$ws.Range("M34").Value = -1266
$ws.Range("N34").Value = -4927.0908
$ws.Range("H39").Value = 5833.3335
$ws.Range("J39").Value = 5833.3335
$ws.Range("L39").Value = 17500.0005
$ws.Range("N39").Value = -18088.0005
$ws.Range("H55").Value = 1833
$ws.Range("J55").Value = 2080.6
$ws.Range("L55").Value = 6241.799999999999
$ws.Range("N55").Value = -6595.799999999999

# GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H98").Value = 8358.4
$ws.Range("J98").Value = 8358.4
$ws.Range("L98").Value = 8358.4
$ws.Range("N98").Value = -14348.4

# LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2280.125
$ws.Range("I61").Value = 2457
$ws.Range("J61").Value = 1749.5
$ws.Range("K61").Value = 2457
$ws.Range("L61").Value = 1749.5
$ws.Range("M61").Value = -2255
$ws.Range("N61").Value = -2153.5
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524
$ws.Range("H113").Value = 2280.125
$ws.Range("I113").Value = 2457
$ws.Range("J113").Value = 1749.5
$ws.Range("K113").Value = 2457
$ws.Range("L113").Value = 1749.5
$ws.Range("M113").Value = -287
$ws.Range("N113").Value = -6089.5
$ws.Range("H136").Value = 3633.3333
$ws.Range("I136").Value = 2960
$ws.Range("K136").Value = 8880
$ws.Range("M136").Value = -6330

# WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3316.6667
$ws.Range("I132").Value = 2980
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8940
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -6410
$ws.Range("N132").Value = -20060
